$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the species-related data between row 3 and row 4,
# keeping the other (identical) columns untouched.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr3 = "{0}3" -f $col
    $addr4 = "{0}4" -f $col

    $val3 = $ws.Range($addr3).Value()
    $val4 = $ws.Range($addr4).Value()

    $ws.Range($addr3).Value = $val4
    $ws.Range($addr4).Value = $val3
}
